{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Summary of the edit (see commit message \"Small adjustments / manual\"):\n//  - Appends a continuation sentence to the end of the \"MyVAO ist eine\n//    Klasse...\" paragraph (\". In der draw() Methode ... modifiziert werden.\")\n//  - Adds several new Heading2/body paragraph pairs documenting\n//    Shader.h/.cpp, Texture.h/.cpp, VertexShader.glsl, FragmentShader.glsl\n//  - Adds a new \"Probleme\" Heading1 section with its explanatory paragraph.\n// (The `w:proofErr` spell-check tags and paragraph-bookmark splitting seen\n// in the raw XML diff are editor-generated artifacts of the text being\n// retyped/merged; they carry no semantic content and are not something the\n// Word object model lets a script author target directly, so this script\n// focuses on the actual textual/structural content change.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// The paragraph that currently ends the document: \"MyVAO ist eine Klasse ... Grafikkarte\"\nconst items = paragraphs.items;\nconst lastParagraph = items[items.length - 1];\n\n// 1) Continue that paragraph with the two extra sentences about draw().\nlastParagraph.insertText(\n  \". In der draw() Methode werden Transformationen, Texturen, Transparenz und Shader gesetzt mit welchen gezeichnet wird. Die Eigenschaften mit welchen das VAO gezeichnet wird k\u00f6nnen mit den Setter-Methoden der Klasse modifiziert werden.\",\n  Word.InsertLocation.end\n);\n\n// 2) Append the new sections after it, each as a Heading 2 + body paragraph,\n//    finishing with a Heading 1 \"Probleme\" section.\nlet anchor = lastParagraph;\n\nfunction addHeading(text, styleName) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n  anchor.style = styleName;\n  return anchor;\n}\n\nfunction addBody(text) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n  anchor.style = \"Normal\";\n  return anchor;\n}\n\naddHeading(\"Shader.h / Shader.cpp\", \"Heading 2\");\naddBody(\n  \"Eine Containerklasse welche Shader von der Festplatte lesen kann, diese Kompiliert, auf die Grafikkarte l\u00e4dt und die Adresse davon in einer Membervariable speichert.\"\n);\n\naddHeading(\"Texture.h / Texture.cpp\", \"Heading 2\");\naddBody(\n  \"Eine Containerklasse welche .dds Texturen von der Festplatte lesen kann, diese auf die Grafikkarte l\u00e4dt und die Adresse davon in einer Membervariable speichert. Diese Klasse kann derzeit nur .dds Dateien lesen. Sie kann vorhandene Mipmaps dieser auslesen und erstellt diese korrekt. Der Compression-Level der .dds files wird auch ausgelesen und dementsprechend gesetzt. \"\n);\n\naddHeading(\"VertexShader.glsl\", \"Heading 2\");\naddBody(\n  \"Vertexshader welcher f\u00fcr die Transformationen von VBOs und deren Normals zust\u00e4ndig ist. Kann derzeit Normals nicht richtig skalieren sofern diese nicht gleichm\u00e4\u00dfig skaliert werden. \"\n);\n\naddHeading(\"FragmentShader.glsl\", \"Heading 2\");\naddBody(\n  \"Shader zust\u00e4ndig f\u00fcr die Farbe eines Fragments. Die Lichtberechnung nach Phong wird hier durchgef\u00fchrt und mit Farbe und Textur kombiniert.\"\n);\n\naddHeading(\"Probleme\", \"Heading 1\");\naddBody(\n  \"Dies war mein erster Versuch mit Shadern zu programmieren und das hat auch sehr gut geklappt, ein Problem konnte ich jedoch nicht l\u00f6sen, die Transparenz. Zwar werden die Objekte je nach Alpha Wert mehr oder weniger opak gezeichnet doch konnte ich die Reihenfolge in der die einzelnen Dreiecke der VBOs gezeichnet werden nicht mehr beeinflussen nachdem die VBOs als ganze Objekte auf der Grafikkarte liegen. Mein Ansatz war daher die Objekte je nach Entfernung zur Kamera zu ordnen und von hinten nach vorne zu zeichnen. Die Transparenz von Objekten gegen\u00fcber anderen Objekten funktioniert so einwandfrei, jedoch werden die Dreiecke eines Objekts in sich in falscher Reihenfolge gerendert und daher Teilweise ausgeblendet.\"\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $d resolves to the active document ($word.ActiveDocument).\n#\n# Summary of the edit (see commit message \"Small adjustments / manual\"):\n#  - Appends a continuation sentence to the end of the \"MyVAO ist eine\n#    Klasse...\" paragraph (\". In der draw() Methode ... modifiziert werden.\")\n#  - Adds several new Heading2/body paragraph pairs documenting\n#    Shader.h/.cpp, Texture.h/.cpp, VertexShader.glsl, FragmentShader.glsl\n#  - Adds a new \"Probleme\" Heading1 section with its explanatory paragraph.\n# (The `w:proofErr` spell-check tags and paragraph-bookmark splitting seen\n# in the raw XML diff are editor-generated artifacts of the text being\n# retyped/merged; they carry no semantic content and are not something the\n# Word object model lets a script author target directly, so this script\n# focuses on the actual textual/structural content change.)\n\n$d = $word.ActiveDocument\n\n# The paragraph that currently ends the document: \"MyVAO ist eine Klasse ... Grafikkarte\"\n$lastParagraph = $d.Paragraphs.Last\n$r = $lastParagraph.Range\n$r.SetRange($r.End - 1, $r.End - 1)\n$r.InsertAfter(\". In der draw() Methode werden Transformationen, Texturen, Transparenz und Shader gesetzt mit welchen gezeichnet wird. Die Eigenschaften mit welchen das VAO gezeichnet wird k\u00f6nnen mit den Setter-Methoden der Klasse modifiziert werden.\")\n\nfunction Add-Heading($text, $styleName) {\n    $d.Paragraphs.Last.Range.InsertParagraphAfter()\n    $p = $d.Paragraphs.Last\n    $p.Range.Text = $text\n    $p.Range.set_Style($styleName)\n}\n\nfunction Add-Body($text) {\n    $d.Paragraphs.Last.Range.InsertParagraphAfter()\n    $p = $d.Paragraphs.Last\n    $p.Range.Text = $text\n    $p.Range.set_Style(\"Normal\")\n}\n\nAdd-Heading \"Shader.h / Shader.cpp\" \"Heading 2\"\nAdd-Body \"Eine Containerklasse welche Shader von der Festplatte lesen kann, diese Kompiliert, auf die Grafikkarte l\u00e4dt und die Adresse davon in einer Membervariable speichert.\"\n\nAdd-Heading \"Texture.h / Texture.cpp\" \"Heading 2\"\nAdd-Body \"Eine Containerklasse welche .dds Texturen von der Festplatte lesen kann, diese auf die Grafikkarte l\u00e4dt und die Adresse davon in einer Membervariable speichert. Diese Klasse kann derzeit nur .dds Dateien lesen. Sie kann vorhandene Mipmaps dieser auslesen und erstellt diese korrekt. Der Compression-Level der .dds files wird auch ausgelesen und dementsprechend gesetzt. \"\n\nAdd-Heading \"VertexShader.glsl\" \"Heading 2\"\nAdd-Body \"Vertexshader welcher f\u00fcr die Transformationen von VBOs und deren Normals zust\u00e4ndig ist. Kann derzeit Normals nicht richtig skalieren sofern diese nicht gleichm\u00e4\u00dfig skaliert werden. \"\n\nAdd-Heading \"FragmentShader.glsl\" \"Heading 2\"\nAdd-Body \"Shader zust\u00e4ndig f\u00fcr die Farbe eines Fragments. Die Lichtberechnung nach Phong wird hier durchgef\u00fchrt und mit Farbe und Textur kombiniert.\"\n\nAdd-Heading \"Probleme\" \"Heading 1\"\nAdd-Body \"Dies war mein erster Versuch mit Shadern zu programmieren und das hat auch sehr gut geklappt, ein Problem konnte ich jedoch nicht l\u00f6sen, die Transparenz. Zwar werden die Objekte je nach Alpha Wert mehr oder weniger opak gezeichnet doch konnte ich die Reihenfolge in der die einzelnen Dreiecke der VBOs gezeichnet werden nicht mehr beeinflussen nachdem die VBOs als ganze Objekte auf der Grafikkarte liegen. Mein Ansatz war daher die Objekte je nach Entfernung zur Kamera zu ordnen und von hinten nach vorne zu zeichnen. Die Transparenz von Objekten gegen\u00fcber anderen Objekten funktioniert so einwandfrei, jedoch werden die Dreiecke eines Objekts in sich in falscher Reihenfolge gerendert und daher Teilweise ausgeblendet.\"\n"}
